$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - classical-best-embed vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.092
$ws.Range("D2").Value = 0.048
$ws.Range("E2").Value = 0.042
$ws.Range("F2").Value = 0.023
$ws.Range("H2").Value = 0.046
$ws.Range("I2").Value = 0.051
$ws.Range("J2").Value = 0.05

# Row 3 - BERT-base vs. classical-best-tfidf (text unchanged)
$ws.Range("C3").Value = 0.119
$ws.Range("D3").Value = 0.147
$ws.Range("E3").Value = 0.141
$ws.Range("F3").Value = 0.106
$ws.Range("G3").Value = 0.13
$ws.Range("H3").Value = 0.135
$ws.Range("I3").Value = 0.128
$ws.Range("J3").Value = 0.13

# Row 4 - BERT-base vs. classical-best-embed
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.027
$ws.Range("D4").Value = 0.099
$ws.Range("E4").Value = 0.099
$ws.Range("F4").Value = 0.083
$ws.Range("G4").Value = 0.081
$ws.Range("H4").Value = 0.089
$ws.Range("I4").Value = 0.077
$ws.Range("J4").Value = 0.08

# Row 5 - BERT-base-nli vs. classical-best-tfidf (text unchanged)
$ws.Range("B5").Value = 0.277
$ws.Range("C5").Value = 0.201
$ws.Range("D5").Value = 0.167
$ws.Range("E5").Value = 0.157
$ws.Range("F5").Value = 0.12
$ws.Range("G5").Value = 0.114
$ws.Range("H5").Value = 0.116
$ws.Range("I5").Value = 0.161
$ws.Range("J5").Value = 0.146

# Row 6 - BERT-base-nli vs. classical-best-embed
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.277
$ws.Range("C6").Value = 0.109
$ws.Range("D6").Value = 0.119
$ws.Range("E6").Value = 0.115
$ws.Range("F6").Value = 0.097
$ws.Range("G6").Value = 0.065
$ws.Range("H6").Value = 0.07000000000000001
$ws.Range("I6").Value = 0.11
$ws.Range("J6").Value = 0.096

# Row 7 - BERT-base-nli vs. BERT-base (text unchanged)
$ws.Range("B7").Value = 0.277
$ws.Range("C7").Value = 0.082
$ws.Range("D7").Value = 0.02
$ws.Range("E7").Value = 0.016
$ws.Range("F7").Value = 0.014
$ws.Range("G7").Value = -0.016
$ws.Range("H7").Value = -0.019
$ws.Range("I7").Value = 0.033
$ws.Range("J7").Value = 0.016
